$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45859.01041666666
$ws.Cells.Item(2, 2).Value = 402
$ws.Cells.Item(3, 1).Value = 45859.02083333334
$ws.Cells.Item(3, 2).Value = 395
$ws.Cells.Item(4, 1).Value = 45859.03125
$ws.Cells.Item(4, 2).Value = 389
$ws.Cells.Item(5, 1).Value = 45859.04166666666
$ws.Cells.Item(5, 2).Value = 397
$ws.Cells.Item(6, 1).Value = 45859.05208333334
$ws.Cells.Item(6, 2).Value = 406
$ws.Cells.Item(7, 1).Value = 45859.0625
$ws.Cells.Item(7, 2).Value = 426
$ws.Cells.Item(8, 1).Value = 45859.07291666666
$ws.Cells.Item(8, 2).Value = 417
$ws.Cells.Item(9, 1).Value = 45859.08333333334
$ws.Cells.Item(9, 2).Value = 408
$ws.Cells.Item(10, 1).Value = 45859.09375
$ws.Cells.Item(10, 2).Value = 420
$ws.Cells.Item(11, 1).Value = 45859.10416666666
$ws.Cells.Item(11, 2).Value = 420
$ws.Cells.Item(12, 1).Value = 45859.11458333334
$ws.Cells.Item(12, 2).Value = 409
$ws.Cells.Item(13, 1).Value = 45859.125
$ws.Cells.Item(13, 2).Value = 399
$ws.Cells.Item(14, 1).Value = 45859.13541666666
$ws.Cells.Item(14, 2).Value = 416
$ws.Cells.Item(15, 1).Value = 45859.14583333334
$ws.Cells.Item(15, 2).Value = 428
$ws.Cells.Item(16, 1).Value = 45859.15625
$ws.Cells.Item(16, 2).Value = 441
$ws.Cells.Item(17, 1).Value = 45859.16666666666
$ws.Cells.Item(17, 2).Value = 413
$ws.Cells.Item(18, 1).Value = 45859.17708333334
$ws.Cells.Item(18, 2).Value = 410
$ws.Cells.Item(19, 1).Value = 45859.1875
$ws.Cells.Item(19, 2).Value = 407
$ws.Cells.Item(20, 1).Value = 45859.19791666666
$ws.Cells.Item(20, 2).Value = 391
$ws.Cells.Item(21, 1).Value = 45859.20833333334
$ws.Cells.Item(21, 2).Value = 375
$ws.Cells.Item(22, 1).Value = 45859.21875
$ws.Cells.Item(22, 2).Value = 379
$ws.Cells.Item(23, 1).Value = 45859.22916666666
$ws.Cells.Item(23, 2).Value = 359
$ws.Cells.Item(24, 1).Value = 45859.23958333334
$ws.Cells.Item(24, 2).Value = 340
$ws.Cells.Item(25, 1).Value = 45859.25
$ws.Cells.Item(25, 2).Value = 309
$ws.Cells.Item(26, 1).Value = 45859.26041666666
$ws.Cells.Item(26, 2).Value = 277
$ws.Cells.Item(27, 1).Value = 45859.27083333334
$ws.Cells.Item(27, 2).Value = 266
$ws.Cells.Item(28, 1).Value = 45859.28125
$ws.Cells.Item(28, 2).Value = 241
$ws.Cells.Item(29, 1).Value = 45859.29166666666
$ws.Cells.Item(29, 2).Value = 206
$ws.Cells.Item(30, 1).Value = 45859.30208333334
$ws.Cells.Item(30, 2).Value = 162
$ws.Cells.Item(31, 1).Value = 45859.3125
$ws.Cells.Item(31, 2).Value = 125
$ws.Cells.Item(32, 1).Value = 45859.32291666666
$ws.Cells.Item(32, 2).Value = 94
$ws.Cells.Item(33, 1).Value = 45859.33333333334
$ws.Cells.Item(33, 2).Value = 68
$ws.Cells.Item(34, 1).Value = 45859.34375
$ws.Cells.Item(34, 2).Value = 47
$ws.Cells.Item(35, 1).Value = 45859.35416666666
$ws.Cells.Item(35, 2).Value = 35
$ws.Cells.Item(36, 1).Value = 45859.36458333334
$ws.Cells.Item(36, 2).Value = 31
$ws.Cells.Item(37, 1).Value = 45859.375
$ws.Cells.Item(37, 2).Value = 29
$ws.Cells.Item(38, 1).Value = 45859.38541666666
$ws.Cells.Item(38, 2).Value = 30
$ws.Cells.Item(39, 1).Value = 45859.39583333334
$ws.Cells.Item(39, 2).Value = 27
$ws.Cells.Item(40, 1).Value = 45859.40625
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(41, 1).Value = 45859.41666666666
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(42, 1).Value = 45859.42708333334
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(43, 1).Value = 45859.4375
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(44, 1).Value = 45859.44791666666
$ws.Cells.Item(44, 2).Value = 0
$ws.Cells.Item(45, 1).Value = 45859.45833333334
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(46, 1).Value = 45859.46875
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(47, 1).Value = 45859.47916666666
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(48, 1).Value = 45859.48958333334
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(49, 1).Value = 45859.5
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(50, 1).Value = 45859.51041666666
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(51, 1).Value = 45859.52083333334
$ws.Cells.Item(51, 2).Value = 0
$ws.Cells.Item(52, 1).Value = 45859.53125
$ws.Cells.Item(52, 2).Value = 0
$ws.Cells.Item(53, 1).Value = 45859.54166666666
$ws.Cells.Item(53, 2).Value = 0
$ws.Cells.Item(54, 1).Value = 45859.55208333334
$ws.Cells.Item(54, 2).Value = 0
$ws.Cells.Item(55, 1).Value = 45859.5625
$ws.Cells.Item(55, 2).Value = 0
$ws.Cells.Item(56, 1).Value = 45859.57291666666
$ws.Cells.Item(56, 2).Value = 0
$ws.Cells.Item(57, 1).Value = 45859.58333333334
$ws.Cells.Item(57, 2).Value = 0
$ws.Cells.Item(58, 1).Value = 45859.59375
$ws.Cells.Item(58, 2).Value = 0
$ws.Cells.Item(59, 1).Value = 45859.60416666666
$ws.Cells.Item(59, 2).Value = 0
$ws.Cells.Item(60, 1).Value = 45859.61458333334
$ws.Cells.Item(60, 2).Value = 0
$ws.Cells.Item(61, 1).Value = 45859.625
$ws.Cells.Item(61, 2).Value = 0
$ws.Cells.Item(62, 1).Value = 45859.63541666666
$ws.Cells.Item(62, 2).Value = 0
$ws.Cells.Item(63, 1).Value = 45859.64583333334
$ws.Cells.Item(63, 2).Value = 0
$ws.Cells.Item(64, 1).Value = 45859.65625
$ws.Cells.Item(64, 2).Value = 0
$ws.Cells.Item(65, 1).Value = 45859.66666666666
$ws.Cells.Item(65, 2).Value = 0
$ws.Cells.Item(66, 1).Value = 45859.67708333334
$ws.Cells.Item(66, 2).Value = 0
$ws.Cells.Item(67, 1).Value = 45859.6875
$ws.Cells.Item(67, 2).Value = 0
$ws.Cells.Item(68, 1).Value = 45859.69791666666
$ws.Cells.Item(68, 2).Value = 0
$ws.Cells.Item(69, 1).Value = 45859.70833333334
$ws.Cells.Item(69, 2).Value = 0
$ws.Cells.Item(70, 1).Value = 45859.71875
$ws.Cells.Item(70, 2).Value = 0
$ws.Cells.Item(71, 1).Value = 45859.72916666666
$ws.Cells.Item(71, 2).Value = 0
$ws.Cells.Item(72, 1).Value = 45859.73958333334
$ws.Cells.Item(72, 2).Value = 0
$ws.Cells.Item(73, 1).Value = 45859.75
$ws.Cells.Item(73, 2).Value = 0
$ws.Cells.Item(74, 1).Value = 45859.76041666666
$ws.Cells.Item(74, 2).Value = 0
$ws.Cells.Item(75, 1).Value = 45859.77083333334
$ws.Cells.Item(75, 2).Value = 0
$ws.Cells.Item(76, 1).Value = 45859.78125
$ws.Cells.Item(76, 2).Value = 0
$ws.Cells.Item(77, 1).Value = 45859.79166666666
$ws.Cells.Item(77, 2).Value = 0
$ws.Cells.Item(78, 1).Value = 45859.80208333334
$ws.Cells.Item(78, 2).Value = 0
$ws.Cells.Item(79, 1).Value = 45859.8125
$ws.Cells.Item(79, 2).Value = 0
$ws.Cells.Item(80, 1).Value = 45859.82291666666
$ws.Cells.Item(80, 2).Value = 0
$ws.Cells.Item(81, 1).Value = 45859.83333333334
$ws.Cells.Item(81, 2).Value = 0
$ws.Cells.Item(82, 1).Value = 45859.84375
$ws.Cells.Item(82, 2).Value = 0
$ws.Cells.Item(83, 1).Value = 45859.85416666666
$ws.Cells.Item(83, 2).Value = 0
$ws.Cells.Item(84, 1).Value = 45859.86458333334
$ws.Cells.Item(84, 2).Value = 0
$ws.Cells.Item(85, 1).Value = 45859.875
$ws.Cells.Item(85, 2).Value = 0
$ws.Cells.Item(86, 1).Value = 45859.88541666666
$ws.Cells.Item(86, 2).Value = 0
$ws.Cells.Item(87, 1).Value = 45859.89583333334
$ws.Cells.Item(87, 2).Value = 0
$ws.Cells.Item(88, 1).Value = 45859.90625
$ws.Cells.Item(88, 2).Value = 0
$ws.Cells.Item(89, 1).Value = 45859.91666666666
$ws.Cells.Item(89, 2).Value = 0
$ws.Cells.Item(90, 1).Value = 45859.92708333334
$ws.Cells.Item(90, 2).Value = 0
$ws.Cells.Item(91, 1).Value = 45859.9375
$ws.Cells.Item(91, 2).Value = 0
$ws.Cells.Item(92, 1).Value = 45859.94791666666
$ws.Cells.Item(92, 2).Value = 0
$ws.Cells.Item(93, 1).Value = 45859.95833333334
$ws.Cells.Item(93, 2).Value = 0
$ws.Cells.Item(94, 1).Value = 45859.96875
$ws.Cells.Item(94, 2).Value = 0
$ws.Cells.Item(95, 1).Value = 45859.97916666666
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(96, 1).Value = 45859.98958333334
$ws.Cells.Item(96, 2).Value = 0
$ws.Cells.Item(97, 1).Value = 45860
$ws.Cells.Item(97, 2).Value = 0
